$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect it so the cell values below can be updated.
$ws.Unprotect()

# Update the confidential disclosure date in the shared string text (cell A59)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."

# Update performance table values (columns D and E, rows 2-56)
$ws.Cells.Item(2, 4).Value = 0.01471438288020881
$ws.Cells.Item(2, 5).Value = 0.00451009132934943
$ws.Cells.Item(3, 4).Value = 0.05141351488013729
$ws.Cells.Item(3, 5).Value = 0.01106718895193337
$ws.Cells.Item(4, 4).Value = 0.01429113307009669
$ws.Cells.Item(4, 5).Value = 0.01057723315059644
$ws.Cells.Item(5, 4).Value = 0.009491168070545567
$ws.Cells.Item(5, 5).Value = 0.006733134792179074
$ws.Cells.Item(6, 4).Value = 0.01542778739463995
$ws.Cells.Item(6, 5).Value = 0.005958426433746444
$ws.Cells.Item(7, 4).Value = 0.01947397692661771
$ws.Cells.Item(7, 5).Value = 0.0003707548568883112
$ws.Cells.Item(8, 4).Value = 0.005061577378551177
$ws.Cells.Item(8, 5).Value = -0.01822817080943273
$ws.Cells.Item(9, 4).Value = 0.006876334676307602
$ws.Cells.Item(9, 5).Value = -0.001608492842206766
$ws.Cells.Item(10, 4).Value = 0.01418098857465158
$ws.Cells.Item(10, 5).Value = -0.005576679340938084
$ws.Cells.Item(11, 4).Value = 0.00797218791578779
$ws.Cells.Item(11, 5).Value = -0.001695699090488723
$ws.Cells.Item(12, 4).Value = 0.01495753031032113
$ws.Cells.Item(12, 5).Value = -0.02538821789499623
$ws.Cells.Item(13, 4).Value = 0.00321846209772757
$ws.Cells.Item(13, 5).Value = -0.03942532576010682
$ws.Cells.Item(14, 4).Value = 0.006052877837557238
$ws.Cells.Item(14, 5).Value = -0.01340033500837501
$ws.Cells.Item(15, 4).Value = 0.01428163944412946
$ws.Cells.Item(15, 5).Value = -0.01696924324661553
$ws.Cells.Item(16, 4).Value = 0.01057387156273133
$ws.Cells.Item(16, 5).Value = -0.01725941422594157
$ws.Cells.Item(17, 4).Value = 0.02241279182837835
$ws.Cells.Item(17, 5).Value = -0.001665535747332081
$ws.Cells.Item(18, 4).Value = 0.008689525064923804
$ws.Cells.Item(18, 5).Value = -0.0007425007425007069
$ws.Cells.Item(19, 4).Value = 0.01621818552290246
$ws.Cells.Item(19, 5).Value = 0.002485451018428675
$ws.Cells.Item(20, 4).Value = 0.0137697517346145
$ws.Cells.Item(20, 5).Value = 0.008969610424383045
$ws.Cells.Item(21, 4).Value = 0.006666061614431343
$ws.Cells.Item(21, 5).Value = -0.00784447476125516
$ws.Cells.Item(22, 4).Value = 0.01470535010987105
$ws.Cells.Item(22, 5).Value = -0.009945009945009797
$ws.Cells.Item(23, 4).Value = 0.01930754659695907
$ws.Cells.Item(23, 5).Value = -0.001374865735767949
$ws.Cells.Item(24, 4).Value = 0.009807099966406696
$ws.Cells.Item(24, 5).Value = -0.01408194785136729
$ws.Cells.Item(25, 4).Value = 0.0213188127351247
$ws.Cells.Item(25, 5).Value = -0.0005245803357314394
$ws.Cells.Item(26, 4).Value = 0.01311042094737539
$ws.Cells.Item(26, 5).Value = -0.01038385826771637
$ws.Cells.Item(27, 4).Value = 0.02170080060454427
$ws.Cells.Item(27, 5).Value = 0.01266138200619826
$ws.Cells.Item(28, 4).Value = 0.05477730011969342
$ws.Cells.Item(28, 5).Value = 0.02457793482528459
$ws.Cells.Item(29, 4).Value = 0.02002313249470308
$ws.Cells.Item(29, 5).Value = -0.03938832252085256
$ws.Cells.Item(30, 4).Value = 0.03090559298693252
$ws.Cells.Item(30, 5).Value = 0.01351996182599025
$ws.Cells.Item(31, 4).Value = 0.01530329492726377
$ws.Cells.Item(31, 5).Value = 0.01998217203981567
$ws.Cells.Item(32, 4).Value = 0.01314139044567627
$ws.Cells.Item(32, 5).Value = -0.0009772565742713901
$ws.Cells.Item(33, 4).Value = 0.01755085710851681
$ws.Cells.Item(33, 5).Value = 0.01449458028737061
$ws.Cells.Item(34, 4).Value = 0.04479885403024699
$ws.Cells.Item(34, 5).Value = 0.007698954818533554
$ws.Cells.Item(35, 4).Value = 0.01081351649005962
$ws.Cells.Item(35, 5).Value = -0.002386634844868674
$ws.Cells.Item(36, 4).Value = 0.009929165260800966
$ws.Cells.Item(36, 5).Value = -0.0001732801940738105
$ws.Cells.Item(37, 4).Value = 0.01071489338535144
$ws.Cells.Item(37, 5).Value = -0.01290322580645165
$ws.Cells.Item(38, 4).Value = 0.007131280010534546
$ws.Cells.Item(38, 5).Value = -0.0002584981258887398
$ws.Cells.Item(39, 4).Value = 0.01196067832295497
$ws.Cells.Item(39, 5).Value = -0.01194457716196851
$ws.Cells.Item(40, 4).Value = 0.0178370177307137
$ws.Cells.Item(40, 5).Value = -0.0109549023187876
$ws.Cells.Item(41, 4).Value = 0.01711906610740504
$ws.Cells.Item(41, 5).Value = -0.004296528677623868
$ws.Cells.Item(42, 4).Value = 0.03451716493660439
$ws.Cells.Item(42, 5).Value = 0.01236348100082774
$ws.Cells.Item(43, 4).Value = 0.01146719411491048
$ws.Cells.Item(43, 5).Value = -0.004072490327835521
$ws.Cells.Item(44, 4).Value = 0.02237841199841932
$ws.Cells.Item(44, 5).Value = -0.003745318352059934
$ws.Cells.Item(45, 4).Value = 0.01301548468770303
$ws.Cells.Item(45, 5).Value = 0.01089866156787767
$ws.Cells.Item(46, 4).Value = 0.008187284599919244
$ws.Cells.Item(46, 5).Value = -0.01046978959100275
$ws.Cells.Item(47, 4).Value = 0.01220530049108161
$ws.Cells.Item(47, 5).Value = -0.007657453556864535
$ws.Cells.Item(48, 4).Value = 0.01039662648763429
$ws.Cells.Item(48, 5).Value = -0.02577196202026655
$ws.Cells.Item(49, 4).Value = 0.01623428474622554
$ws.Cells.Item(49, 5).Value = -0.004405786569164016
$ws.Cells.Item(50, 4).Value = 0.008575693724442866
$ws.Cells.Item(50, 5).Value = -0.02204762774834035
$ws.Cells.Item(51, 4).Value = 0.01035923573422935
$ws.Cells.Item(51, 5).Value = -0.008452609038656633
$ws.Cells.Item(52, 4).Value = 0.00820028072867051
$ws.Cells.Item(52, 5).Value = 0.0003297065611607142
$ws.Cells.Item(53, 4).Value = 0.009015134929003962
$ws.Cells.Item(53, 5).Value = -0.00968557660193714
$ws.Cells.Item(54, 4).Value = 0.1334567852911105
$ws.Cells.Item(55, 4).Value = 0.04429129836364913
$ws.Cells.Item(55, 5).Value = 0.001533018867924563
$ws.Cells.Item(56, 5).Value = 0.0005265736422408462

# Restore sheet protection
$ws.Protect()

Write-Host "Applied all updates"
